# Apply timekeeper export fixes:
#  - revert admin dev default (Rate/Total zeroed out back to 0 on seeded rows)
#  - re-seed customer names with the real (non-dev) customer list
#  - refresh the employee id to the new seeded value
#  - clear the leftover "Seeded sample hours" note (now blank)

$wb = $excel.ActiveWorkbook

$wsTime = $wb.Worksheets.Item("Weekly Timesheet")
$wsSchema = $wb.Worksheets.Item("Jason Schema")

# New customer names for rows 2-6 (Date order 01-12 .. 01-16)
$newClients = @("Ueltschi", "Landers", "Jackson / Ho", "Bottomley", "Lynn")
$newEmpId = "emp_fn0y5dge"

# --- Weekly Timesheet sheet ---
for ($i = 0; $i -lt 5; $i++) {
    $row = 2 + $i
    $wsTime.Range("B$row").Value = $newClients[$i]
    $wsTime.Range("E$row").Value = 0
    $wsTime.Range("F$row").Value = 0
}

# Subtotal / hourly subtotal / grand total rows revert to 0
$wsTime.Range("F8").Value = 0
$wsTime.Range("F11").Value = 0
$wsTime.Range("F13").Value = 0

# --- Jason Schema sheet ---
for ($i = 0; $i -lt 5; $i++) {
    $row = 2 + $i
    $wsSchema.Range("B$row").Value = $newEmpId
    $wsSchema.Range("D$row").Value = $newClients[$i]
    $wsSchema.Range("F$row").Value = 0
    $wsSchema.Range("G$row").Value = 0
    $wsSchema.Range("I$row").Value = ""
}

$wb.Save()
